$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Sheet ALC row 76 (G=12602)
$ws.Cells.Item(76, 8).Value = 10305.375
$ws.Cells.Item(76, 10).Value = 7315
$ws.Cells.Item(76, 12).Value = 7315
$ws.Cells.Item(76, 14).Value = -7945

# Sheet ALC row 79 (G=12602)
$ws.Cells.Item(79, 8).Value = 10305.375
$ws.Cells.Item(79, 10).Value = 7315
$ws.Cells.Item(79, 12).Value = 7315
$ws.Cells.Item(79, 14).Value = -9499

# Sheet ALC row 86 (G=12603)
$ws.Cells.Item(86, 8).Value = 53498.75
$ws.Cells.Item(86, 9).Value = 5000
$ws.Cells.Item(86, 10).Value = 101997.5
$ws.Cells.Item(86, 11).Value = 5000
$ws.Cells.Item(86, 12).Value = 101997.5
$ws.Cells.Item(86, 13).Value = -3877
$ws.Cells.Item(86, 14).Value = -104243.5

# Sheet ALC row 89 (G=12603)
$ws.Cells.Item(89, 8).Value = 53498.75
$ws.Cells.Item(89, 9).Value = 5000
$ws.Cells.Item(89, 10).Value = 101997.5
$ws.Cells.Item(89, 11).Value = 25000
$ws.Cells.Item(89, 12).Value = 509987.5
$ws.Cells.Item(89, 13).Value = -19384
$ws.Cells.Item(89, 14).Value = -521219.5

# Sheet ALC row 107 (G=27766)
$ws.Cells.Item(107, 8).Value = 20002284
$ws.Cells.Item(107, 9).Value = 23810770
$ws.Cells.Item(107, 10).Value = 7724.5
$ws.Cells.Item(107, 11).Value = 23810770
$ws.Cells.Item(107, 12).Value = 7724.5
$ws.Cells.Item(107, 13).Value = -23808850
$ws.Cells.Item(107, 14).Value = -11564.5

# Sheet ALC row 132 (G=44049)
$ws.Cells.Item(132, 8).Value = 7698.579
$ws.Cells.Item(132, 9).Value = 3976.1936
$ws.Cells.Item(132, 11).Value = 11928.5808
$ws.Cells.Item(132, 13).Value = -9398.5808

# Sheet ALC row 135 (G=44047)
$ws.Cells.Item(135, 8).Value = 1007.5
$ws.Cells.Item(135, 9).Value = 648.1177
$ws.Cells.Item(135, 10).Value = 2229.4
$ws.Cells.Item(135, 11).Value = 5833.0593
$ws.Cells.Item(135, 12).Value = 20064.6
$ws.Cells.Item(135, 13).Value = -3298.0593
$ws.Cells.Item(135, 14).Value = -25134.6

$ws = $wb.Worksheets.Item("ARM")
# Sheet ARM row 2 (G=27713)
$ws.Cells.Item(2, 8).Value = 45456710
$ws.Cells.Item(2, 9).Value = 62502070
$ws.Cells.Item(2, 11).Value = 62502070
$ws.Cells.Item(2, 13).Value = -62501957

# Sheet ARM row 37 (G=3096)
$ws.Cells.Item(37, 8).Value = 0
$ws.Cells.Item(37, 10).Value = 0
$ws.Cells.Item(37, 12).ClearContents()
$ws.Cells.Item(37, 14).Value = 0

# Sheet ARM row 61 (G=43999)
$ws.Cells.Item(61, 8).Value = 16067
$ws.Cells.Item(61, 9).Value = 15575.375
$ws.Cells.Item(61, 10).Value = 20000
$ws.Cells.Item(61, 11).Value = 15575.375
$ws.Cells.Item(61, 12).Value = 20000
$ws.Cells.Item(61, 13).Value = -15363.375
$ws.Cells.Item(61, 14).Value = -20424

# Sheet ARM row 74 (G=44000)
$ws.Cells.Item(74, 8).Value = 7745.5
$ws.Cells.Item(74, 9).Value = 7407.143
$ws.Cells.Item(74, 11).Value = 7407.143
$ws.Cells.Item(74, 13).Value = -6533.143

# Sheet ARM row 77 (G=44000)
$ws.Cells.Item(77, 8).Value = 7745.5
$ws.Cells.Item(77, 9).Value = 7407.143
$ws.Cells.Item(77, 11).Value = 37035.715
$ws.Cells.Item(77, 13).Value = -32667.715

# Sheet ARM row 116 (G=27713)
$ws.Cells.Item(116, 8).Value = 45456710
$ws.Cells.Item(116, 9).Value = 62502070
$ws.Cells.Item(116, 11).Value = 62502070
$ws.Cells.Item(116, 13).Value = -62499776

# Sheet ARM row 132 (G=43997)
$ws.Cells.Item(132, 8).Value = 3453.3333
$ws.Cells.Item(132, 9).Value = 3454.634
$ws.Cells.Item(132, 11).Value = 10363.902
$ws.Cells.Item(132, 13).Value = -7833.902

# Sheet ARM row 136 (G=43999)
$ws.Cells.Item(136, 8).Value = 16067
$ws.Cells.Item(136, 9).Value = 15575.375
$ws.Cells.Item(136, 10).Value = 20000
$ws.Cells.Item(136, 11).Value = 46726.125
$ws.Cells.Item(136, 12).Value = 60000
$ws.Cells.Item(136, 13).Value = -44176.125
$ws.Cells.Item(136, 14).Value = -65100

$ws = $wb.Worksheets.Item("BSM")
# Sheet BSM row 3 (G=27713)
$ws.Cells.Item(3, 8).Value = 45456710
$ws.Cells.Item(3, 9).Value = 62502070
$ws.Cells.Item(3, 11).Value = 62502070
$ws.Cells.Item(3, 13).Value = -62501956

# Sheet BSM row 22 (G=5092)
$ws.Cells.Item(22, 8).Value = 562.1
$ws.Cells.Item(22, 9).Value = 562.1
$ws.Cells.Item(22, 10).Value = 0
$ws.Cells.Item(22, 11).Value = 562.1
$ws.Cells.Item(22, 12).Value = 0
$ws.Cells.Item(22, 13).ClearContents()
$ws.Cells.Item(22, 14).Value = -389.1

# Sheet BSM row 107 (G=27706)
$ws.Cells.Item(107, 8).Value = 9228.666999999999
$ws.Cells.Item(107, 9).Value = 9399.5
$ws.Cells.Item(107, 10).Value = 8887
$ws.Cells.Item(107, 11).Value = 9399.5
$ws.Cells.Item(107, 12).Value = 8887
$ws.Cells.Item(107, 13).Value = -7479.5
$ws.Cells.Item(107, 14).Value = -12727

# Sheet BSM row 134 (G=43998)
$ws.Cells.Item(134, 8).Value = 6275.607
$ws.Cells.Item(134, 9).Value = 6529.12
$ws.Cells.Item(134, 10).Value = 4163
$ws.Cells.Item(134, 11).Value = 19587.36
$ws.Cells.Item(134, 12).Value = 12489
$ws.Cells.Item(134, 13).Value = -17052.36
$ws.Cells.Item(134, 14).Value = -17559

$ws = $wb.Worksheets.Item("CRP")
# Sheet CRP row 11 (G=1821)
$ws.Cells.Item(11, 8).Value = 1299
$ws.Cells.Item(11, 9).Value = 999
$ws.Cells.Item(11, 11).Value = 999
$ws.Cells.Item(11, 13).Value = -859

# Sheet CRP row 132 (G=44019)
$ws.Cells.Item(132, 8).Value = 5013.5835
$ws.Cells.Item(132, 9).Value = 5002.8774
$ws.Cells.Item(132, 10).Value = 5061.273
$ws.Cells.Item(132, 11).Value = 15008.6322
$ws.Cells.Item(132, 12).Value = 15183.819
$ws.Cells.Item(132, 13).Value = -12478.6322
$ws.Cells.Item(132, 14).Value = -20243.819

# Sheet CRP row 134 (G=44020)
$ws.Cells.Item(134, 8).Value = 4298.019
$ws.Cells.Item(134, 9).Value = 4679.15
$ws.Cells.Item(134, 10).Value = 3027.5833
$ws.Cells.Item(134, 11).Value = 14037.45
$ws.Cells.Item(134, 12).Value = 9082.749899999999
$ws.Cells.Item(134, 13).Value = -11502.45
$ws.Cells.Item(134, 14).Value = -14152.7499

$ws = $wb.Worksheets.Item("CUL")
# Sheet CUL row 2 (G=4847)
$ws.Cells.Item(2, 8).Value = 139.25
$ws.Cells.Item(2, 9).Value = 71.333336
$ws.Cells.Item(2, 11).Value = 428.000016
$ws.Cells.Item(2, 13).Value = -315.000016

# Sheet CUL row 3 (G=44094)
$ws.Cells.Item(3, 8).Value = 3015.1333
$ws.Cells.Item(3, 9).Value = 1940.5385
$ws.Cells.Item(3, 10).Value = 10000
$ws.Cells.Item(3, 11).Value = 5821.6155
$ws.Cells.Item(3, 12).Value = 30000
$ws.Cells.Item(3, 13).Value = -5709.6155
$ws.Cells.Item(3, 14).Value = -30224

# Sheet CUL row 50 (G=4725)
$ws.Cells.Item(50, 8).Value = 2351.6428
$ws.Cells.Item(50, 10).Value = 3119.3
$ws.Cells.Item(50, 12).Value = 9357.900000000001
$ws.Cells.Item(50, 14).Value = -10319.9

# Sheet CUL row 53 (G=4725)
$ws.Cells.Item(53, 8).Value = 2351.6428
$ws.Cells.Item(53, 10).Value = 3119.3
$ws.Cells.Item(53, 12).Value = 9357.900000000001
$ws.Cells.Item(53, 14).Value = -10319.9

# Sheet CUL row 54 (G=4671)
$ws.Cells.Item(54, 8).Value = 599
$ws.Cells.Item(54, 10).Value = 599
$ws.Cells.Item(54, 12).Value = 1797
$ws.Cells.Item(54, 14).Value = -2915

# Sheet CUL row 57 (G=4655)
$ws.Cells.Item(57, 8).Value = 6187.125
$ws.Cells.Item(57, 10).Value = 5213.857
$ws.Cells.Item(57, 12).Value = 15641.571
$ws.Cells.Item(57, 14).Value = -16759.571

# Sheet CUL row 60 (G=4750)
$ws.Cells.Item(60, 8).Value = 170
$ws.Cells.Item(60, 9).Value = 182.33333
$ws.Cells.Item(60, 10).Value = 151.5
$ws.Cells.Item(60, 11).Value = 546.99999
$ws.Cells.Item(60, 12).Value = 454.5
$ws.Cells.Item(60, 13).Value = -295.99999
$ws.Cells.Item(60, 14).Value = -956.5

# Sheet CUL row 69 (G=12850)
$ws.Cells.Item(69, 8).Value = 942
$ws.Cells.Item(69, 10).Value = 0
$ws.Cells.Item(69, 12).Value = 0
$ws.Cells.Item(69, 14).ClearContents()

# Sheet CUL row 72 (G=12850)
$ws.Cells.Item(72, 8).Value = 942
$ws.Cells.Item(72, 10).Value = 0
$ws.Cells.Item(72, 12).Value = 0
$ws.Cells.Item(72, 14).ClearContents()

# Sheet CUL row 134 (G=44074)
$ws.Cells.Item(134, 8).Value = 1967.4286
$ws.Cells.Item(134, 9).Value = 1967.4286
$ws.Cells.Item(134, 10).Value = 0
$ws.Cells.Item(134, 11).Value = 5902.2858
$ws.Cells.Item(134, 12).Value = 0
$ws.Cells.Item(134, 13).ClearContents()
$ws.Cells.Item(134, 14).Value = -832.2857999999997

$ws = $wb.Worksheets.Item("GSM")
# Sheet GSM row 80 (G=12521)
$ws.Cells.Item(80, 8).Value = 88753000
$ws.Cells.Item(80, 9).Value = 118335464
$ws.Cells.Item(80, 11).Value = 118335464
$ws.Cells.Item(80, 13).Value = -118334466

# Sheet GSM row 83 (G=12521)
$ws.Cells.Item(83, 8).Value = 88753000
$ws.Cells.Item(83, 9).Value = 118335464
$ws.Cells.Item(83, 11).Value = 591677320
$ws.Cells.Item(83, 13).Value = -591672328

# Sheet GSM row 102 (G=36169)
$ws.Cells.Item(102, 8).Value = 6348.2666
$ws.Cells.Item(102, 9).Value = 6709.5386
$ws.Cells.Item(102, 10).Value = 4000
$ws.Cells.Item(102, 11).Value = 6709.5386
$ws.Cells.Item(102, 12).Value = 4000
$ws.Cells.Item(102, 13).Value = -5087.5386
$ws.Cells.Item(102, 14).Value = -7244

# Sheet GSM row 132 (G=44008)
$ws.Cells.Item(132, 8).Value = 3884.4
$ws.Cells.Item(132, 9).Value = 3605.75
$ws.Cells.Item(132, 10).Value = 4999
$ws.Cells.Item(132, 11).Value = 10817.25
$ws.Cells.Item(132, 12).Value = 14997
$ws.Cells.Item(132, 13).Value = -8287.25
$ws.Cells.Item(132, 14).Value = -20057

$ws = $wb.Worksheets.Item("LTW")
# Sheet LTW row 7 (G=36249)
$ws.Cells.Item(7, 8).Value = 3543.25
$ws.Cells.Item(7, 9).Value = 2164.25
$ws.Cells.Item(7, 10).Value = 4922.25
$ws.Cells.Item(7, 11).Value = 2164.25
$ws.Cells.Item(7, 12).Value = 4922.25
$ws.Cells.Item(7, 13).Value = -2052.25
$ws.Cells.Item(7, 14).Value = -5146.25

# Sheet LTW row 55 (G=5284)
$ws.Cells.Item(55, 8).Value = 1487.825
$ws.Cells.Item(55, 9).Value = 1321.3103
$ws.Cells.Item(55, 11).Value = 1321.3103
$ws.Cells.Item(55, 13).Value = -1148.3103

# Sheet LTW row 126 (G=36249)
$ws.Cells.Item(126, 8).Value = 3543.25
$ws.Cells.Item(126, 9).Value = 2164.25
$ws.Cells.Item(126, 10).Value = 4922.25
$ws.Cells.Item(126, 11).Value = 6492.75
$ws.Cells.Item(126, 12).Value = 14766.75
$ws.Cells.Item(126, 13).Value = -4022.75
$ws.Cells.Item(126, 14).Value = -19706.75

$ws = $wb.Worksheets.Item("WVR")
# Sheet WVR row 34 (G=3349)
$ws.Cells.Item(34, 8).Value = 0
$ws.Cells.Item(34, 9).Value = 0
$ws.Cells.Item(34, 11).Value = 0
$ws.Cells.Item(34, 13).ClearContents()

# Sheet WVR row 38 (G=27990)
$ws.Cells.Item(38, 8).Value = 24963.334
$ws.Cells.Item(38, 9).Value = 24900
$ws.Cells.Item(38, 10).Value = 24995
$ws.Cells.Item(38, 11).Value = 24900
$ws.Cells.Item(38, 12).Value = 24995
$ws.Cells.Item(38, 13).Value = -24427
$ws.Cells.Item(38, 14).Value = -25941

# Sheet WVR row 52 (G=2816)
$ws.Cells.Item(52, 8).Value = 12999.667
$ws.Cells.Item(52, 9).Value = 5000
$ws.Cells.Item(52, 10).Value = 16999.5
$ws.Cells.Item(52, 11).Value = 5000
$ws.Cells.Item(52, 12).Value = 16999.5
$ws.Cells.Item(52, 13).Value = -4774
$ws.Cells.Item(52, 14).Value = -17451.5

# Sheet WVR row 53 (G=3172)
$ws.Cells.Item(53, 8).Value = 39994.5
$ws.Cells.Item(53, 9).Value = 39994.5
$ws.Cells.Item(53, 11).Value = 39994.5
$ws.Cells.Item(53, 13).Value = -39387.5

# Sheet WVR row 55 (G=2832)
$ws.Cells.Item(55, 8).Value = 8023.5
$ws.Cells.Item(55, 9).Value = 8023.5
$ws.Cells.Item(55, 11).Value = 8023.5
$ws.Cells.Item(55, 13).Value = -7746.5

# Sheet WVR row 107 (G=27746)
$ws.Cells.Item(107, 8).Value = 1960.5333
$ws.Cells.Item(107, 9).Value = 2615.3
$ws.Cells.Item(107, 10).Value = 651
$ws.Cells.Item(107, 11).Value = 7845.900000000001
$ws.Cells.Item(107, 12).Value = 1953
$ws.Cells.Item(107, 13).Value = -5925.900000000001
$ws.Cells.Item(107, 14).Value = -5793
